$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '51.761.82'
$ws.Range("E2").Value = '  +0.40%  '

# Row 3
$ws.Range("D3").Value = '3.114.24'
$ws.Range("E3").Value = '  +4.12%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '388.44'
$ws.Range("E5").Value = '  +1.65%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.41'
$ws.Range("E6").Value = '  -0.88%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  -0.79%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.47'
$ws.Range("E10").Value = '  +1.73%  '

# Row 11
$ws.Range("E11").Value = '  +0.03%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0862'
$ws.Range("E12").Value = '  +0.02%  '

# Row 13
$ws.Range("D13").Value = '3.605.07'
$ws.Range("E13").Value = '  +3.89%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.93'
$ws.Range("E14").Value = '  +1.79%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.72'
$ws.Range("E15").Value = '  +1.05%  '

# Row 16
$ws.Range("D16").Value = '3.116.48'
$ws.Range("E16").Value = '  +4.57%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.998'
$ws.Range("E17").Value = '  +0.27%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.04'
$ws.Range("E18").Value = '  -2.04%  '

# Row 19
$ws.Range("D19").Value = '51.881.06'
$ws.Range("E19").Value = '  +0.48%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.20'
$ws.Range("E20").Value = '  +3.84%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.48'
$ws.Range("E21").Value = '  -0.65%  '

# Row 22
$ws.Range("E22").Value = '  +0.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.09'
$ws.Range("E23").Value = '  -0.46%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.91'
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("E25").Value = '  -3.01%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.16'
$ws.Range("E26").Value = '  +1.46%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.18'
$ws.Range("E27").Value = '  +4.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.171'
$ws.Range("E28").Value = '  +0.78%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.21'
$ws.Range("E29").Value = '  +0.07%  '

# Row 30
$ws.Range("E30").Value = '  -0.01%  '

# Row 31
$ws.Range("E31").Value = '  -0.79%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.38'
$ws.Range("E32").Value = '  -0.43%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.77'
$ws.Range("E33").Value = '  +3.00%  '

# Row 34
$ws.Range("E34").Value = '  +0.74%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.30'
$ws.Range("E35").Value = '  -2.15%  '

# Row 36
$ws.Range("E36").Value = '  +0.95%  '

# Row 37
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.40'
$ws.Range("E38").Value = '  +3.29%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.289'
$ws.Range("E39").Value = '  +6.18%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("E40").Value = '  +2.72%  '

# Row 41
$ws.Range("E41").Value = '  +1.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '129.39'
$ws.Range("E42").Value = '  +1.12%  '

# Row 43
$ws.Range("E43").Value = '  -1.46%  '

# Row 44
$ws.Range("E44").Value = '  -0.09%  '

# Row 45
$ws.Range("E45").Value = '  -4.08%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.37'
$ws.Range("E46").Value = '  +4.20%  '

# Row 47
$ws.Range("E47").Value = '  +5.15%  '

# Row 48
$ws.Range("E48").Value = '  +2.14%  '

# Row 49
$ws.Range("D49").Value = '2.052.46'
$ws.Range("E49").Value = '  +0.75%  '

# Row 50
$ws.Range("D50").Value = '3.418.56'
$ws.Range("E50").Value = '  +4.03%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.208'
$ws.Range("E51").Value = '  +5.32%  '
